# TaskStatus.xlsx maintenance update
# - For each data row (2..55) on Sheet1:
#     * Column B (Period)         -> new period value
#     * Column C (LastSuccessful) -> cleared (task not yet completed)
#     * Column D (NextDue)        -> refreshed due date
#     * Column E (active)         -> flipped from 1 to 0
# - Refreshes the frozen-pane scroll position and current selection
#   to reflect where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row, NewPeriod, NewNextDue(serial), NewActive
$rows = @(
    @(2, 1, 45254, 0),
    @(3, 2, 45254, 0),
    @(4, 3, 45254, 0),
    @(5, 5, 45254, 0),
    @(6, 7, 45254, 0),
    @(7, 30, 45254, 0),
    @(8, 1, 45254, 0),
    @(9, 2, 45254, 0),
    @(10, 3, 45254, 0),
    @(11, 5, 45254, 0),
    @(12, 7, 45254, 0),
    @(13, 7, 45254, 0),
    @(14, 30, 45254, 0),
    @(15, 1, 45254, 0),
    @(16, 1, 45254, 0),
    @(17, 7, 45254, 0),
    @(18, 2, 45254, 0),
    @(19, 3, 45254, 0),
    @(20, 5, 45254, 0),
    @(21, 7, 45254, 0),
    @(22, 1, 45254, 0),
    @(23, 7, 45254, 0),
    @(24, 2, 45254, 0),
    @(25, 3, 45254, 0),
    @(26, 5, 45254, 0),
    @(27, 1, 45254, 0),
    @(28, 7, 45254, 0),
    @(29, 2, 45254, 0),
    @(30, 3, 45254, 0),
    @(31, 7, 45254, 0),
    @(32, 1, 45254, 0),
    @(33, 2, 45254, 0),
    @(34, 3, 45254, 0),
    @(35, 7, 45254, 0),
    @(36, 7, 45254, 0),
    @(37, 1, 45254, 0),
    @(38, 2, 45254, 0),
    @(39, 3, 45254, 0),
    @(40, 30, 45254, 0),
    @(41, 1, 45254, 0),
    @(42, 7, 45254, 0),
    @(43, 1, 45254, 0),
    @(44, 2, 45254, 0),
    @(45, 3, 45255, 0),
    @(46, 1, 45255, 0),
    @(47, 7, 45255, 0),
    @(48, 1, 45255, 0),
    @(49, 2, 45256, 0),
    @(50, 3, 45256, 0),
    @(51, 1, 45256, 0),
    @(52, 7, 45256, 0),
    @(53, 1, 45257, 0),
    @(54, 2, 45257, 0),
    @(55, 3, 45257, 0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $period = $r[1]
    $nextDue = $r[2]
    $active = $r[3]

    $ws.Cells.Item($rowNum, 2).Value = $period
    $ws.Cells.Item($rowNum, 3).ClearContents()
    $ws.Cells.Item($rowNum, 4).Value = $nextDue
    $ws.Cells.Item($rowNum, 5).Value = $active
}

# Restore the view: frozen pane scrolled back to the top, selection on
# the last edited block of rows.
$ws.Range("B51:B55").Select()
